# The "holiday_ad_content" column (E) contains text that was built by
# concatenating the Campaign, Date Range and Advertisement fields with a
# real TAB character and a real newline as separators. The fix replaces
# those real control characters with the literal two-character escape
# sequences "\t" and "\n" (backslash + letter), leaving everything else
# in the string unchanged.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 61
$col = 5  # column E

for ($row = $firstRow; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, $col)
    $text = $cell.Value2
    if ($text -ne $null) {
        $text = $text.Replace([char]9, "\t")
        $text = $text.Replace([char]10, "\n")
        $cell.Value = $text
    }
}
